$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 4.889
$ws.Range("F2").Value = 0.669
$ws.Range("D3").Value = 4.889
$ws.Range("F3").Value = 1.311
$ws.Range("D4").Value = 4.889
$ws.Range("F4").Value = 0.839
$ws.Range("D5").Value = 4.889
$ws.Range("F5").Value = 0.934
$ws.Range("D6").Value = 4.889
$ws.Range("F6").Value = 1.248
$ws.Range("D7").Value = 6.023
$ws.Range("F7").Value = 1.313
$ws.Range("D8").Value = 6.023
$ws.Range("F8").Value = 0.959
$ws.Range("D9").Value = 6.023
$ws.Range("F9").Value = 0.872
$ws.Range("D10").Value = 6.023
$ws.Range("F10").Value = 0.665
$ws.Range("D11").Value = 6.023
$ws.Range("F11").Value = 1.192
$ws.Range("D12").Value = 4.798
$ws.Range("F12").Value = 1.247
$ws.Range("D13").Value = 4.798
$ws.Range("F13").Value = 1.047
$ws.Range("D14").Value = 4.798
$ws.Range("F14").Value = 0.688
$ws.Range("D15").Value = 4.798
$ws.Range("F15").Value = 1.019
$ws.Range("D16").Value = 3.196
$ws.Range("F16").Value = 0.802
$ws.Range("D17").Value = 3.196
$ws.Range("F17").Value = 0.886
$ws.Range("D18").Value = 3.196
$ws.Range("F18").Value = 1.312
$ws.Range("D19").Value = 3.523
$ws.Range("F19").Value = 1.321
$ws.Range("D20").Value = 3.523
$ws.Range("F20").Value = 0.679
$ws.Range("D21").Value = 1.414
$ws.Range("F21").Value = 0.916
$ws.Range("D22").Value = 1.414
$ws.Range("F22").Value = 1.397
$ws.Range("D23").Value = 1.414
$ws.Range("F23").Value = 0.972
$ws.Range("D24").Value = 1.414
$ws.Range("F24").Value = 0.715
$ws.Range("D25").Value = 2.111
$ws.Range("F25").Value = 0.938
$ws.Range("D26").Value = 2.111
$ws.Range("F26").Value = 0.777
$ws.Range("D27").Value = 2.111
$ws.Range("F27").Value = 1.053
$ws.Range("D28").Value = 2.111
$ws.Range("F28").Value = 0.612
$ws.Range("D29").Value = 2.111
$ws.Range("F29").Value = 1.619
$ws.Range("D30").Value = 1.718
$ws.Range("F30").Value = 1.22
$ws.Range("D31").Value = 1.718
$ws.Range("F31").Value = 1.091
$ws.Range("D32").Value = 1.718
$ws.Range("F32").Value = 0.797
$ws.Range("D33").Value = 1.718
$ws.Range("F33").Value = 0.892
$ws.Range("D34").Value = 3.967
$ws.Range("F34").Value = 0.632
$ws.Range("D35").Value = 3.967
$ws.Range("F35").Value = 1.368
$ws.Range("D36").Value = 4.961
$ws.Range("F36").Value = 0.956
$ws.Range("D37").Value = 4.961
$ws.Range("F37").Value = 1.034
$ws.Range("D38").Value = 4.961
$ws.Range("F38").Value = 1.01
$ws.Range("D39").Value = 4.864
$ws.Range("F39").Value = 1.015
$ws.Range("D40").Value = 4.864
$ws.Range("F40").Value = 1.112
$ws.Range("D41").Value = 4.864
$ws.Range("F41").Value = 0.703
$ws.Range("D42").Value = 4.864
$ws.Range("F42").Value = 1.169
